$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

foreach ($row in 2,3) {
    $ws.Range("D$row").Value = -0.534
    $ws.Range("G$row").Value = -0.34375
    $ws.Range("H$row").Value = -0.34375
    $ws.Range("I$row").Value = 0.0625
    $ws.Range("J$row").Value = 0.03125
    $ws.Range("K$row").Value = -0.001
    $ws.Range("L$row").Value = -0.03125
    $ws.Range("U$row").Value = 0.001
    $ws.Range("V$row").Value = 0.002857142857142858
    $ws.Range("W$row").Value = 0.009009009009009009
    $ws.Range("X$row").Value = 0.08213270210227022
    $ws.Range("Y$row").Value = -0.07312369309326121
    $ws.Range("Z$row").Value = -0.2831858407079646
    $ws.Range("AA$row").Value = -0.008849557522123894
    $ws.Range("AB$row").Value = 0.08213270210227022
    $ws.Range("AC$row").Value = -0.09098225962439412
    $ws.Range("AG$row").Value = -0.001
    $ws.Range("AJ$row").Value = -0.002865329512893983
    $ws.Range("AK$row").Value = 0.009433962264150943
    $ws.Range("AN$row").Value = 0
    $ws.Range("AP$row").Value = -0.5
}
